# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the Cilantro data block
# (rows 237-238), pushing the existing rows 237-244 down to 239-246.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 237:244 down to 239:246 by inserting two new rows.
$ws.Rows("237:238").Insert()

# New row 237: "Primera" quality entry for the latest week.
$ws.Range("A237").Value = 7
$ws.Range("B237").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C237").Value = 'Ñuble'
$ws.Range("D237").Value = 45075
$ws.Range("E237").Value = 16
$ws.Range("F237").Value = 100112040
$ws.Range("G237").Value = 'Cilantro'
$ws.Range("H237").Value = 'Sin especificar'
$ws.Range("I237").Value = 'Primera'
$ws.Range("J237").Value = 150
$ws.Range("K237").Value = 1200
$ws.Range("L237").Value = 1200
$ws.Range("M237").Value = 1200
$ws.Range("N237").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O237").Value = 'Provincia de Diguillín'
$ws.Range("P237").Value = 1200
$ws.Range("Q237").Value = 1
$ws.Range("R237").Value = 'Hortaliza'

# New row 238: "Segunda" quality entry for the latest week.
$ws.Range("A238").Value = 7
$ws.Range("B238").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C238").Value = 'Ñuble'
$ws.Range("D238").Value = 45075
$ws.Range("E238").Value = 16
$ws.Range("F238").Value = 100112040
$ws.Range("G238").Value = 'Cilantro'
$ws.Range("H238").Value = 'Sin especificar'
$ws.Range("I238").Value = 'Segunda'
$ws.Range("J238").Value = 100
$ws.Range("K238").Value = 1000
$ws.Range("L238").Value = 1000
$ws.Range("M238").Value = 1000
$ws.Range("N238").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O238").Value = 'Provincia de Diguillín'
$ws.Range("P238").Value = 1000
$ws.Range("Q238").Value = 1
$ws.Range("R238").Value = 'Hortaliza'

Write-Output "Inserted rows 237-238; sheet now spans to row 246"
